$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows (2-16) got their weekly records (Fecha + Volumen +
# Precio minimo/maximo/promedio ponderado + Precio $/Kg) re-shuffled across
# rows. Apply the new values for each affected row/column directly.

$updates = @{
    2  = @{ D = 44504 }
    3  = @{ D = 44516 }
    4  = @{ D = 44523; J = 400 }
    5  = @{ D = 44524; J = 400; K = 800;  L = 900;  M = 850;  P = 850 }
    6  = @{ D = 44511; J = 500; K = 900;  L = 1000; M = 950;  P = 950 }
    7  = @{ D = 44517; J = 500; K = 800;  L = 900;  M = 850;  P = 850 }
    8  = @{ D = 44530; J = 300; K = 800;  L = 900;  M = 850;  P = 850 }
    9  = @{ D = 44510; J = 600 }
    10 = @{ D = 44476; J = 300; K = 1100; L = 1200; M = 1150; P = 1150 }
    11 = @{ D = 44503; K = 900;  L = 1000; M = 950;  P = 950 }
    12 = @{ D = 44505; J = 440; K = 900;  L = 1000; M = 950;  P = 950 }
    13 = @{ D = 44518; J = 400; K = 800;  L = 900;  M = 850;  P = 850 }
    14 = @{ D = 44525; J = 360; K = 800;  L = 900;  M = 850;  P = 850 }
    15 = @{ D = 44508; J = 400; K = 900;  L = 1000; M = 950;  P = 950 }
    16 = @{ D = 44512; J = 600 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
